# Fill in the newly-marked XMOS package pin columns (B, F, J, N) on the
# "XMOS Dualchip" sheet. These columns already had header labels in row 1
# ("U1 Package Pin" / "U2 Package Pin") but the per-row pin values were
# blank; the commit adds the package-pin text for rows 2-45.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("XMOS Dualchip")

$ws.Range("B2").Value = "A15"
$ws.Range("F2").Value = "A39"
$ws.Range("J2").Value = "A15"
$ws.Range("N2").Value = "A39"
$ws.Range("B3").Value = "A16"
$ws.Range("F3").Value = "A40"
$ws.Range("J3").Value = "A16"
$ws.Range("N3").Value = "A40"
$ws.Range("B4").Value = "A5"
$ws.Range("F4").Value = "A41"
$ws.Range("J4").Value = "A5"
$ws.Range("N4").Value = "A41"
$ws.Range("B5").Value = "A6"
$ws.Range("F5").Value = "A42"
$ws.Range("J5").Value = "A6"
$ws.Range("N5").Value = "A42"
$ws.Range("B6").Value = "A7"
$ws.Range("F6").Value = "A43"
$ws.Range("J6").Value = "A7"
$ws.Range("N6").Value = "A43"
$ws.Range("B7").Value = "A8"
$ws.Range("F7").Value = "A44"
$ws.Range("J7").Value = "A8"
$ws.Range("N7").Value = "A44"
$ws.Range("B8").Value = "A9"
$ws.Range("F8").Value = "A45"
$ws.Range("J8").Value = "A9"
$ws.Range("N8").Value = "A45"
$ws.Range("B9").Value = "A10"
$ws.Range("F9").Value = "A46"
$ws.Range("J9").Value = "A10"
$ws.Range("N9").Value = "A46"
$ws.Range("B10").Value = "A11"
$ws.Range("F10").Value = "A47"
$ws.Range("J10").Value = "A11"
$ws.Range("N10").Value = "A47"
$ws.Range("B11").Value = "A12"
$ws.Range("F11").Value = "A48"
$ws.Range("J11").Value = "A12"
$ws.Range("N11").Value = "A48"
$ws.Range("B12").Value = "A13"
$ws.Range("F12").Value = "A49"
$ws.Range("J12").Value = "A13"
$ws.Range("N12").Value = "A49"
$ws.Range("B13").Value = "A14"
$ws.Range("F13").Value = "A50"
$ws.Range("J13").Value = "A14"
$ws.Range("N13").Value = "A50"
$ws.Range("B14").Value = "A55"
$ws.Range("F14").Value = "A21"
$ws.Range("J14").Value = "A55"
$ws.Range("N14").Value = "A21"
$ws.Range("B15").Value = "A57"
$ws.Range("F15").Value = "A23"
$ws.Range("J15").Value = "A57"
$ws.Range("N15").Value = "A23"
$ws.Range("B16").Value = "A58"
$ws.Range("F16").Value = "A24"
$ws.Range("J16").Value = "A58"
$ws.Range("N16").Value = "A24"
$ws.Range("B17").Value = "A59"
$ws.Range("F17").Value = "A25"
$ws.Range("J17").Value = "A59"
$ws.Range("N17").Value = "A25"
$ws.Range("B18").Value = "A60"
$ws.Range("F18").Value = "A26"
$ws.Range("J18").Value = "A60"
$ws.Range("N18").Value = "A26"
$ws.Range("B19").Value = "A61"
$ws.Range("F19").Value = "A27"
$ws.Range("J19").Value = "A61"
$ws.Range("N19").Value = "A27"
$ws.Range("B20").Value = "A63"
$ws.Range("F20").Value = "A29"
$ws.Range("J20").Value = "A63"
$ws.Range("N20").Value = "A29"
$ws.Range("B21").Value = "A64"
$ws.Range("F21").Value = "A30"
$ws.Range("J21").Value = "A64"
$ws.Range("N21").Value = "A30"
$ws.Range("B22").Value = "A65"
$ws.Range("F22").Value = "A31"
$ws.Range("J22").Value = "A65"
$ws.Range("N22").Value = "A31"
$ws.Range("B23").Value = "A66"
$ws.Range("F23").Value = "A32"
$ws.Range("J23").Value = "A66"
$ws.Range("N23").Value = "A32"
$ws.Range("B24").Value = "A56"
$ws.Range("F24").Value = "A22"
$ws.Range("J24").Value = "A56"
$ws.Range("N24").Value = "A22"
$ws.Range("B25").Value = "A62"
$ws.Range("F25").Value = "A28"
$ws.Range("J25").Value = "A62"
$ws.Range("N25").Value = "A28"
$ws.Range("B26").Value = "A54"
$ws.Range("F26").Value = "A20"
$ws.Range("J26").Value = "A54"
$ws.Range("N26").Value = "A20"
$ws.Range("B27").Value = "A67"
$ws.Range("F27").Value = "A33"
$ws.Range("J27").Value = "A67"
$ws.Range("N27").Value = "A33"
$ws.Range("B28").Value = "B38"
$ws.Range("F28").Value = "B16"
$ws.Range("J28").Value = "B38"
$ws.Range("N28").Value = "B16"
$ws.Range("B29").Value = "B39"
$ws.Range("F29").Value = "B17"
$ws.Range("J29").Value = "B39"
$ws.Range("N29").Value = "B17"
$ws.Range("B30").Value = "B40"
$ws.Range("F30").Value = "B18"
$ws.Range("J30").Value = "B40"
$ws.Range("N30").Value = "B18"
$ws.Range("B31").Value = "B41"
$ws.Range("F31").Value = "B19"
$ws.Range("J31").Value = "B41"
$ws.Range("N31").Value = "B19"
$ws.Range("B32").Value = "B44"
$ws.Range("F32").Value = "B20"
$ws.Range("J32").Value = "B44"
$ws.Range("N32").Value = "B20"
$ws.Range("B33").Value = "B45"
$ws.Range("F33").Value = "B21"
$ws.Range("J33").Value = "B45"
$ws.Range("N33").Value = "B21"
$ws.Range("B34").Value = "B46"
$ws.Range("F34").Value = "B22"
$ws.Range("J34").Value = "B46"
$ws.Range("N34").Value = "B22"
$ws.Range("B35").Value = "B47"
$ws.Range("F35").Value = "B23"
$ws.Range("J35").Value = "B47"
$ws.Range("N35").Value = "B23"
$ws.Range("B36").Value = "A4"
$ws.Range("F36").Value = "B24"
$ws.Range("J36").Value = "A4"
$ws.Range("N36").Value = "B24"
$ws.Range("B37").Value = "A3"
$ws.Range("F37").Value = "B25"
$ws.Range("J37").Value = "A3"
$ws.Range("N37").Value = "B25"
$ws.Range("B38").Value = "B48"
$ws.Range("F38").Value = "B26"
$ws.Range("J38").Value = "B48"
$ws.Range("N38").Value = "B26"
$ws.Range("B39").Value = "B49"
$ws.Range("F39").Value = "B27"
$ws.Range("J39").Value = "B49"
$ws.Range("N39").Value = "B27"
$ws.Range("B40").Value = "B50"
$ws.Range("F40").Value = "B30"
$ws.Range("J40").Value = "B50"
$ws.Range("N40").Value = "B30"
$ws.Range("B41").Value = "B51"
$ws.Range("F41").Value = "B31"
$ws.Range("J41").Value = "B51"
$ws.Range("N41").Value = "B31"
$ws.Range("B42").Value = "B52"
$ws.Range("J42").Value = "B52"
$ws.Range("B43").Value = "B53"
$ws.Range("J43").Value = "B53"
$ws.Range("B44").Value = "B54"
$ws.Range("J44").Value = "B54"
$ws.Range("B45").Value = "B55"
$ws.Range("J45").Value = "B55"

# The edit also moves the active tab from "XMOS Dualchip Planning" to
# "XMOS Dualchip", and updates that sheet's scroll/selection so the
# newly-filled rows are visible.
$ws.Activate()
$ws.Range("J40").Select()
